$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11 (old row 11 shifts down to row 12,
# preserving all of its existing values/formatting).
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Cells.Item(11, 1).Value = 3
$ws.Cells.Item(11, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44606
$ws.Cells.Item(11, 5).Value = 5
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100101
$ws.Cells.Item(11, 8).Value = "Berries"
$ws.Cells.Item(11, 9).Value = 100101004
$ws.Cells.Item(11, 10).Value = "Frambuesa"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 45
$ws.Cells.Item(11, 14).Value = 7000
$ws.Cells.Item(11, 15).Value = 7000
$ws.Cells.Item(11, 16).Value = 7000
$ws.Cells.Item(11, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia de Linares"
$ws.Cells.Item(11, 19).Value = 3500
$ws.Cells.Item(11, 20).Value = 2
